# Finished Week 13 logging
$wb = $excel.ActiveWorkbook

# "OFF" sheet (sheet1.xml) - row 3 (label "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 395
$wsOff.Range("C3").Value = 295
$wsOff.Range("D3").Value = 91
$wsOff.Range("E3").Value = 43
$wsOff.Range("F3").Value = 7

# "DEF" sheet (sheet2.xml) - row 3 (label "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 385
$wsDef.Range("C3").Value = 250
$wsDef.Range("D3").Value = 101
$wsDef.Range("E3").Value = 52
$wsDef.Range("F3").Value = 6
